# Update colour error name in modifycoin

$wb = $excel.ActiveWorkbook

# --- "Sheet1": numeric-looking text values (exchange-rate style strings) ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# A1 "0,160977" -> "0,224235"
# NOTE: assigning this string straight to .Value gets auto-coerced to a
# number by the locale-aware input parser (comma read as thousands
# separator). Route it through a text formula + paste-as-values so the
# literal text is preserved without touching the cell's style.
$ws1.Range("A1").Formula = '="0,224235"'
$ws1.Range("A1").Copy()
$ws1.Range("A1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# A2 "1866,02" -> "29860,91" (not ambiguous, plain assignment is fine)
$ws1.Range("A2").Value = "29860,91"

# --- "data" sheet: coin names / sheet references ---
$ws2 = $wb.Worksheets.Item("data")

$ws2.Range("B1").Value = "-"
$ws2.Range("C1").Value = "LITECOIN"
$ws2.Range("D1").Value = "MOVER"

$ws2.Range("B2").Value = ""
$ws2.Range("C2").Value = "Sheet1"

$ws2.Range("A3").Value = "A1"
$ws2.Range("B3").Value = ""
$ws2.Range("C3").Value = "A1"
